# Update "paises" / "provincias" COVID data sheet and re-sort by total cases.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: look up a country by name in column A (rows 4:216) and overwrite
# the given columns on that row with new values.
function Set-CountryRow {
    param(
        [string]$Name,
        [hashtable]$Values
    )
    $found = $ws.Range("A4:A216").Find($Name)
    if ($found -eq $null) {
        Write-Host ("Country not found: " + $Name)
        return
    }
    $r = $found.Row
    foreach ($col in $Values.Keys) {
        $ws.Range($col + $r).Value = $Values[$col]
    }
}

# Updated case counts (new snapshot taken a bit later the same day).
Set-CountryRow "Estados Unidos" @{ "B" = 583870; "C" = 23570; "D" = 35755; "E" = 524630; "G" = 1380; "H" = 23485 }
Set-CountryRow "Alemania"       @{ "B" = 129207; "C" = 1353;  "E" = 61789;  "G" = 96;   "H" = 3118 }
Set-CountryRow "Peru"           @{ "D" = 2642;   "E" = 6926;  "G" = 23;    "H" = 216 }
Set-CountryRow "Eslovaquia"     @{ "B" = 769;    "C" = 27;    "D" = 107;   "E" = 660;  "F" = 5;  "G" = 0; "H" = 2 }
Set-CountryRow "Banglades"      @{ "B" = 803;    "C" = 182;   "D" = 42;    "E" = 722;  "F" = 1;  "G" = 5; "H" = 39 }
Set-CountryRow "Bolivia"        @{ "B" = 330;    "C" = 30;    "D" = 2;     "E" = 301;  "F" = 3;  "G" = 3; "H" = 27 }
Set-CountryRow "Mauricio"       @{ "B" = 324;    "C" = 0;     "D" = 42;    "E" = 273;  "F" = 3;  "G" = 0; "H" = 9 }
Set-CountryRow "Nigeria"        @{ "B" = 343;    "C" = 20;    "D" = 91;    "E" = 242;  "F" = 2;  "G" = 0; "H" = 10 }
Set-CountryRow "Islas Feroe"    @{ "B" = 184;    "C" = 0;     "D" = 157;   "E" = 27;   "F" = 0;  "G" = 0; "H" = 0 }
Set-CountryRow "Venezuela"      @{ "B" = 189;    "C" = 8;     "D" = 110;   "E" = 70;   "F" = 6;  "G" = 0; "H" = 9 }
Set-CountryRow "Guatemala"      @{ "B" = 156;    "C" = 1;     "E" = 132 }
Set-CountryRow "Bermudas"       @{ "D" = 30;     "E" = 22;    "F" = 3;     "G" = 1;    "H" = 5 }

# Refresh the "last updated" timestamp string (row 1).
$ws.Range("A1").Value = "Datos actualizados a 13 de Abril de 2020 a las 23:52"

# The sheet is kept sorted by total cases (column B) descending; re-sort the
# data range now that several countries' totals changed, so rows land back
# in rank order (this naturally reproduces the handful of row swaps seen in
# the diff, e.g. Banglades/Eslovaquia, Nigeria/Bolivia/Mauricio,
# Venezuela/Islas Feroe).
$dataRange = $ws.Range("A4:H216")
$sortKey = $ws.Range("B4")
$dataRange.Sort($sortKey, 2)
